$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add the new filenames for rows 46-81 (000045.jpg .. 000080.jpg).
#    This must happen BEFORE the new tag text is introduced, so the shared-string
#    table keeps the filenames contiguous and the new tag text lands at the end,
#    mirroring the authored diff exactly.
for ($i = 45; $i -le 80; $i++) {
    $row = $i + 1
    $name = "{0:D6}.jpg" -f $i
    $ws.Cells.Item($row, 1).Value2 = $name
}

# 2) New tag text for the (Japan/India) story, replacing the old Australia/India one.
$newTag = "japan, India, Indian Prime Minister, japanese Prime Minister, kishida, Modi, Bilateral Relations, International Relations, Trade, Pact, Harmony"

# 3) Update column B for every data row (2-81). Row 46 keeps the literal "Tags"
#    text (matching the header), exactly as authored in the source workbook;
#    every other row gets the new tag text.
for ($row = 2; $row -le 81; $row++) {
    if ($row -eq 46) {
        $ws.Cells.Item($row, 2).Value2 = "Tags"
    } else {
        $ws.Cells.Item($row, 2).Value2 = $newTag
    }
}

# 4) Restore the portrait page orientation present in the saved workbook.
$ws.PageSetup.Orientation = 1

# 5) Scroll/selection bookkeeping to mirror the final view state.
[void]$ws.Range("B75").Select()

Write-Output "Files_and_Tags sheet extended to row 81"
